# Frontend_VideoTopics.pptx — slide 12 ("9. Typical List component composition")
# Content Placeholder 2 text updates:
#   Paragraph 1: merge the two runs of "(See random example ... not exact match)"
#                into a single run (text itself is unchanged).
#   Paragraph 2: change "The example is about Buildings, e.g. three buldings on
#                one campus, two on other." to "The example is a list of Buildings"
#                (split across 3 runs: "The example is a ", "list of ", "Buildings").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

$shape = $null
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "Content Placeholder 2") {
        $shape = $sh
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$tf = $shape.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: collapse the two existing runs into one run ----------------
$para1 = $tr.Paragraphs(1, 1)
# Force a real text-diff so the engine re-merges the whole paragraph into a
# single run (writing the exact same text is treated as a no-op otherwise).
$para1.Text = "placeholder"
$para1b = $tr.Paragraphs(1, 1)
$para1b.Text = "(See random example on the next slide. Idea important, not exact match)"

# --- Paragraph 2: rewrite text, split into 3 runs -----------------------------
$para2 = $tr.Paragraphs(2, 1)
$para2Start = $para2.Start

# Reset the whole paragraph first (also a real text diff) so stale formatting
# (e.g. the err="1" spell-flag on the old "buldings" run) doesn't leak into
# the freshly written runs below.
$para2.Text = "placeholder text goes here"
$para2b = $tr.Paragraphs(2, 1)
$para2b.Text = "The example is a list of Buildings"

$run1 = $tr.Characters($para2Start, 17)
$run1.Text = "The example is a "

$run2 = $tr.Characters($para2Start + 17, 8)
$run2.Text = "list of "

$run3 = $tr.Characters($para2Start + 25, 9)
$run3.Text = "Buildings"
